# Replicate the prior "Criminal Offenses - Non Forcible Sex Offenses" append step,
# but for Offense = "Criminal Offenses - Incest", with Date values sum2014 and sum2015
# only (9 sectors x 4 reporting locations x 2 dates = 72 new rows: 758-829).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$offenseText = "Criminal Offenses - Incest"

$newRows = @(
  [PSCustomObject]@{ Row=758; A="Public, 4-year or above"; B="On Campus (excluding Residence Halls)"; D="sum2014"; E="0" },
  [PSCustomObject]@{ Row=759; A="Private nonprofit, 4-year or above"; B="On Campus (excluding Residence Halls)"; D="sum2014"; E="0" },
  [PSCustomObject]@{ Row=760; A="Private for-profit, 4-year or above"; B="On Campus (excluding Residence Halls)"; D="sum2014"; E="0" },
  [PSCustomObject]@{ Row=761; A="Public, 2-year"; B="On Campus (excluding Residence Halls)"; D="sum2014"; E="0" },
  [PSCustomObject]@{ Row=762; A="Private nonprofit, 2-year"; B="On Campus (excluding Residence Halls)"; D="sum2014"; E="0" },
  [PSCustomObject]@{ Row=763; A="Private for-profit, 2-year"; B="On Campus (excluding Residence Halls)"; D="sum2014"; E="0" },
  [PSCustomObject]@{ Row=764; A="Public, less-than 2-year"; B="On Campus (excluding Residence Halls)"; D="sum2014"; E="0" },
  [PSCustomObject]@{ Row=765; A="Private nonprofit, less-than 2-year"; B="On Campus (excluding Residence Halls)"; D="sum2014"; E="0" },
  [PSCustomObject]@{ Row=766; A="Private for-profit, less-than 2-year"; B="On Campus (excluding Residence Halls)"; D="sum2014"; E="0" },
  [PSCustomObject]@{ Row=767; A="Public, 4-year or above"; B="On Campus (excluding Residence Halls)"; D="sum2015"; E="0" },
  [PSCustomObject]@{ Row=768; A="Private nonprofit, 4-year or above"; B="On Campus (excluding Residence Halls)"; D="sum2015"; E="0" },
  [PSCustomObject]@{ Row=769; A="Private for-profit, 4-year or above"; B="On Campus (excluding Residence Halls)"; D="sum2015"; E="0" },
  [PSCustomObject]@{ Row=770; A="Public, 2-year"; B="On Campus (excluding Residence Halls)"; D="sum2015"; E="0" },
  [PSCustomObject]@{ Row=771; A="Private nonprofit, 2-year"; B="On Campus (excluding Residence Halls)"; D="sum2015"; E="0" },
  [PSCustomObject]@{ Row=772; A="Private for-profit, 2-year"; B="On Campus (excluding Residence Halls)"; D="sum2015"; E="0" },
  [PSCustomObject]@{ Row=773; A="Public, less-than 2-year"; B="On Campus (excluding Residence Halls)"; D="sum2015"; E="0" },
  [PSCustomObject]@{ Row=774; A="Private nonprofit, less-than 2-year"; B="On Campus (excluding Residence Halls)"; D="sum2015"; E="0" },
  [PSCustomObject]@{ Row=775; A="Private for-profit, less-than 2-year"; B="On Campus (excluding Residence Halls)"; D="sum2015"; E="0" },
  [PSCustomObject]@{ Row=776; A="Public, 4-year or above"; B="On Campus (Residence Halls)"; D="sum2014"; E="1" },
  [PSCustomObject]@{ Row=777; A="Private nonprofit, 4-year or above"; B="On Campus (Residence Halls)"; D="sum2014"; E="0" },
  [PSCustomObject]@{ Row=778; A="Private for-profit, 4-year or above"; B="On Campus (Residence Halls)"; D="sum2014"; E="0" },
  [PSCustomObject]@{ Row=779; A="Public, 2-year"; B="On Campus (Residence Halls)"; D="sum2014"; E="0" },
  [PSCustomObject]@{ Row=780; A="Private nonprofit, 2-year"; B="On Campus (Residence Halls)"; D="sum2014"; E="0" },
  [PSCustomObject]@{ Row=781; A="Private for-profit, 2-year"; B="On Campus (Residence Halls)"; D="sum2014"; E="0" },
  [PSCustomObject]@{ Row=782; A="Public, less-than 2-year"; B="On Campus (Residence Halls)"; D="sum2014"; E="0" },
  [PSCustomObject]@{ Row=783; A="Private nonprofit, less-than 2-year"; B="On Campus (Residence Halls)"; D="sum2014"; E="NULL" },
  [PSCustomObject]@{ Row=784; A="Private for-profit, less-than 2-year"; B="On Campus (Residence Halls)"; D="sum2014"; E="0" },
  [PSCustomObject]@{ Row=785; A="Public, 4-year or above"; B="On Campus (Residence Halls)"; D="sum2015"; E="1" },
  [PSCustomObject]@{ Row=786; A="Private nonprofit, 4-year or above"; B="On Campus (Residence Halls)"; D="sum2015"; E="0" },
  [PSCustomObject]@{ Row=787; A="Private for-profit, 4-year or above"; B="On Campus (Residence Halls)"; D="sum2015"; E="0" },
  [PSCustomObject]@{ Row=788; A="Public, 2-year"; B="On Campus (Residence Halls)"; D="sum2015"; E="0" },
  [PSCustomObject]@{ Row=789; A="Private nonprofit, 2-year"; B="On Campus (Residence Halls)"; D="sum2015"; E="0" },
  [PSCustomObject]@{ Row=790; A="Private for-profit, 2-year"; B="On Campus (Residence Halls)"; D="sum2015"; E="0" },
  [PSCustomObject]@{ Row=791; A="Public, less-than 2-year"; B="On Campus (Residence Halls)"; D="sum2015"; E="0" },
  [PSCustomObject]@{ Row=792; A="Private nonprofit, less-than 2-year"; B="On Campus (Residence Halls)"; D="sum2015"; E="NULL" },
  [PSCustomObject]@{ Row=793; A="Private for-profit, less-than 2-year"; B="On Campus (Residence Halls)"; D="sum2015"; E="0" },
  [PSCustomObject]@{ Row=794; A="Public, 4-year or above"; B="Non-Campus"; D="sum2014"; E="0" },
  [PSCustomObject]@{ Row=795; A="Private nonprofit, 4-year or above"; B="Non-Campus"; D="sum2014"; E="0" },
  [PSCustomObject]@{ Row=796; A="Private for-profit, 4-year or above"; B="Non-Campus"; D="sum2014"; E="0" },
  [PSCustomObject]@{ Row=797; A="Public, 2-year"; B="Non-Campus"; D="sum2014"; E="0" },
  [PSCustomObject]@{ Row=798; A="Private nonprofit, 2-year"; B="Non-Campus"; D="sum2014"; E="0" },
  [PSCustomObject]@{ Row=799; A="Private for-profit, 2-year"; B="Non-Campus"; D="sum2014"; E="0" },
  [PSCustomObject]@{ Row=800; A="Public, less-than 2-year"; B="Non-Campus"; D="sum2014"; E="0" },
  [PSCustomObject]@{ Row=801; A="Private nonprofit, less-than 2-year"; B="Non-Campus"; D="sum2014"; E="0" },
  [PSCustomObject]@{ Row=802; A="Private for-profit, less-than 2-year"; B="Non-Campus"; D="sum2014"; E="0" },
  [PSCustomObject]@{ Row=803; A="Public, 4-year or above"; B="Non-Campus"; D="sum2015"; E="0" },
  [PSCustomObject]@{ Row=804; A="Private nonprofit, 4-year or above"; B="Non-Campus"; D="sum2015"; E="0" },
  [PSCustomObject]@{ Row=805; A="Private for-profit, 4-year or above"; B="Non-Campus"; D="sum2015"; E="0" },
  [PSCustomObject]@{ Row=806; A="Public, 2-year"; B="Non-Campus"; D="sum2015"; E="0" },
  [PSCustomObject]@{ Row=807; A="Private nonprofit, 2-year"; B="Non-Campus"; D="sum2015"; E="0" },
  [PSCustomObject]@{ Row=808; A="Private for-profit, 2-year"; B="Non-Campus"; D="sum2015"; E="0" },
  [PSCustomObject]@{ Row=809; A="Public, less-than 2-year"; B="Non-Campus"; D="sum2015"; E="0" },
  [PSCustomObject]@{ Row=810; A="Private nonprofit, less-than 2-year"; B="Non-Campus"; D="sum2015"; E="0" },
  [PSCustomObject]@{ Row=811; A="Private for-profit, less-than 2-year"; B="Non-Campus"; D="sum2015"; E="0" },
  [PSCustomObject]@{ Row=812; A="Public, 4-year or above"; B="Public Property"; D="sum2014"; E="0" },
  [PSCustomObject]@{ Row=813; A="Private nonprofit, 4-year or above"; B="Public Property"; D="sum2014"; E="0" },
  [PSCustomObject]@{ Row=814; A="Private for-profit, 4-year or above"; B="Public Property"; D="sum2014"; E="0" },
  [PSCustomObject]@{ Row=815; A="Public, 2-year"; B="Public Property"; D="sum2014"; E="0" },
  [PSCustomObject]@{ Row=816; A="Private nonprofit, 2-year"; B="Public Property"; D="sum2014"; E="0" },
  [PSCustomObject]@{ Row=817; A="Private for-profit, 2-year"; B="Public Property"; D="sum2014"; E="0" },
  [PSCustomObject]@{ Row=818; A="Public, less-than 2-year"; B="Public Property"; D="sum2014"; E="0" },
  [PSCustomObject]@{ Row=819; A="Private nonprofit, less-than 2-year"; B="Public Property"; D="sum2014"; E="0" },
  [PSCustomObject]@{ Row=820; A="Private for-profit, less-than 2-year"; B="Public Property"; D="sum2014"; E="0" },
  [PSCustomObject]@{ Row=821; A="Public, 4-year or above"; B="Public Property"; D="sum2015"; E="0" },
  [PSCustomObject]@{ Row=822; A="Private nonprofit, 4-year or above"; B="Public Property"; D="sum2015"; E="0" },
  [PSCustomObject]@{ Row=823; A="Private for-profit, 4-year or above"; B="Public Property"; D="sum2015"; E="0" },
  [PSCustomObject]@{ Row=824; A="Public, 2-year"; B="Public Property"; D="sum2015"; E="0" },
  [PSCustomObject]@{ Row=825; A="Private nonprofit, 2-year"; B="Public Property"; D="sum2015"; E="0" },
  [PSCustomObject]@{ Row=826; A="Private for-profit, 2-year"; B="Public Property"; D="sum2015"; E="0" },
  [PSCustomObject]@{ Row=827; A="Public, less-than 2-year"; B="Public Property"; D="sum2015"; E="0" },
  [PSCustomObject]@{ Row=828; A="Private nonprofit, less-than 2-year"; B="Public Property"; D="sum2015"; E="0" },
  [PSCustomObject]@{ Row=829; A="Private for-profit, less-than 2-year"; B="Public Property"; D="sum2015"; E="0" }
)

# Row 757 is the last existing data row; A/B/D/E there use the workbook's
# established "data row" style (cellXfs index carrying quotePrefix). Use
# copy/paste-special (formats only) from it as a template for each new row
# so the new cells reuse that existing style instead of allocating new ones.
$templateRange = $ws.Range("A757:E757")

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $offenseText
    $ws.Cells.Item($row, 4).Value = $r.D
    if ($r.E -eq "NULL") {
        $ws.Cells.Item($row, 5).Value = 0
    } else {
        $ws.Cells.Item($row, 5).Value = [double]$r.E
    }

    $templateRange.Copy()
    $ws.Range("A" + $row + ":E" + $row).PasteSpecial(-4122)

    if ($r.E -eq "NULL") {
        $ws.Cells.Item($row, 5).ClearContents()
    }
}

# Row 776 (sector "Public, 4-year or above" / "On Campus (Residence Halls)" / sum2014)
# carries a stray, empty, styled F cell in the source data (mirroring the same
# pattern seen at the start of every other "On Campus (Residence Halls)" block
# for that sector, e.g. row 731). Reproduce it the same way: copy a template
# that includes that F cell (format-only), which leaves F776 styled but blank.
$ws.Range("A731:F731").Copy()
$ws.Range("A776:F776").PasteSpecial(-4122)

# Move the viewport/selection to mirror the edit location.
$ws.Range("A830").Select()

Write-Host "Appended rows 758-829 for Criminal Offenses - Incest."
